# Apply the "usage_2025_06" report edit:
#  - Header row 3 gains a new "Kampus dan Ruangan" column (inserted before the
#    existing "Tanggal Pemasangan" column, which moves from C to D).
#  - The single remaining data row (row 4) is rewritten to describe a Router
#    installed in room D112 on 2025-06-05.
#  - The other former data rows (5-8) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3): insert "Kampus dan Ruangan" as the new column C,
#     pushing "Tanggal Pemasangan" out to column D. ---
$ws.Range("D3").Value = "Tanggal Pemasangan"
$ws.Range("C3").Value = "Kampus dan Ruangan"

# --- Data row (row 4): now Nama Inventaris / Kampus dan Ruangan / Tanggal Pemasangan ---
$ws.Range("A4").Value = "Router"
$ws.Range("B4").Value = "D112"

# Force the cell to text *before* assigning, so the "2025-06-05" string is
# stored verbatim (shared string) rather than being reinterpreted as a date
# serial number.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2025-06-05"

# --- Remove the old extra data rows (former rows 5-8) ---
$ws.Range("A5:D8").EntireRow.Delete()
